# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice* / LevePrice* / LeveProfit* columns H:N) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with newly pulled price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 608.6667
$ws.Range("I61").Value = 608.6667
$ws.Range("K61").Value = 1826.0001
$ws.Range("M61").Value = -1654.0001
$ws.Range("H69").Value = 5996.3335
$ws.Range("I69").Value = 5994.5
$ws.Range("J69").Value = 6000
$ws.Range("K69").Value = 17983.5
$ws.Range("L69").Value = 18000
$ws.Range("M69").Value = -17109.5
$ws.Range("N69").Value = -19748
$ws.Range("H72").Value = 5996.3335
$ws.Range("I72").Value = 5994.5
$ws.Range("J72").Value = 6000
$ws.Range("K72").Value = 53950.5
$ws.Range("L72").Value = 54000
$ws.Range("M72").Value = -49582.5
$ws.Range("N72").Value = -62736
$ws.Range("H138").Value = 6153.7744
$ws.Range("J138").Value = 10905.875
$ws.Range("L138").Value = 32717.625
$ws.Range("N138").Value = -42997.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1814059.4
$ws.Range("I32").Value = 1955492.2
$ws.Range("J32").Value = 3718.2
$ws.Range("K32").Value = 1955492.2
$ws.Range("L32").Value = 3718.2
$ws.Range("M32").Value = -1955205.2
$ws.Range("N32").Value = -4292.2
$ws.Range("H45").Value = 7205.143
$ws.Range("J45").Value = 10858.875
$ws.Range("L45").Value = 10858.875
$ws.Range("N45").Value = -11612.875
$ws.Range("H132").Value = 7119.189
$ws.Range("I132").Value = 5975.625
$ws.Range("J132").Value = 9230.385
$ws.Range("K132").Value = 17926.875
$ws.Range("L132").Value = 27691.155
$ws.Range("M132").Value = -15396.875
$ws.Range("N132").Value = -32751.155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1664.36
$ws.Range("I105").Value = 1124.7646
$ws.Range("J105").Value = 2811
$ws.Range("K105").Value = 1124.7646
$ws.Range("L105").Value = 2811
$ws.Range("M105").Value = 622.2354
$ws.Range("N105").Value = -6305
$ws.Range("H134").Value = 5581.6665
$ws.Range("I134").Value = 2301.3044
$ws.Range("K134").Value = 6903.9132
$ws.Range("M134").Value = -4368.9132

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 8036.643
$ws.Range("I134").Value = 3771.6924
$ws.Range("K134").Value = 11315.0772
$ws.Range("M134").Value = -8780.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 633446.25
$ws.Range("I4").Value = 2000649.8
$ws.Range("J4").Value = 11990.091
$ws.Range("K4").Value = 6001949.4
$ws.Range("L4").Value = 35970.273
$ws.Range("M4").Value = -6001837.4
$ws.Range("N4").Value = -36194.273
$ws.Range("H5").Value = 4448083.5
$ws.Range("I5").Value = 8000751
$ws.Range("J5").Value = 7249.25
$ws.Range("K5").Value = 24002253
$ws.Range("L5").Value = 21747.75
$ws.Range("M5").Value = -24002141
$ws.Range("N5").Value = -21971.75
$ws.Range("H133").Value = 2309.3333
$ws.Range("I133").Value = 964
$ws.Range("K133").Value = 2892
$ws.Range("M133").Value = 2168
$ws.Range("H135").Value = 4448083.5
$ws.Range("I135").Value = 8000751
$ws.Range("J135").Value = 7249.25
$ws.Range("K135").Value = 72006759
$ws.Range("L135").Value = 65243.25
$ws.Range("M135").Value = -72004224
$ws.Range("N135").Value = -70313.25
$ws.Range("H139").Value = 32243
$ws.Range("I139").Value = 60970.176
$ws.Range("K139").Value = 182910.528
$ws.Range("M139").Value = -177770.528

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5988.829
$ws.Range("I113").Value = 2855.1052
$ws.Range("J113").Value = 8695.227999999999
$ws.Range("K113").Value = 2855.1052
$ws.Range("L113").Value = 8695.227999999999
$ws.Range("M113").Value = -685.1052
$ws.Range("N113").Value = -13035.228
$ws.Range("H132").Value = 7428.8667
$ws.Range("I132").Value = 2608.5715
$ws.Range("K132").Value = 7825.7145
$ws.Range("M132").Value = -5295.7145
$ws.Range("H136").Value = 57170.332
$ws.Range("J136").Value = 57170.332
$ws.Range("L136").Value = 171510.996
$ws.Range("N136").Value = -176610.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4099.1665
$ws.Range("I16").Value = 3289.6365
$ws.Range("K16").Value = 3289.6365
$ws.Range("M16").Value = -3119.6365
$ws.Range("H61").Value = 2634204
$ws.Range("I61").Value = 3847095.2
$ws.Range("J61").Value = 6272.8335
$ws.Range("K61").Value = 3847095.2
$ws.Range("L61").Value = 6272.8335
$ws.Range("M61").Value = -3846893.2
$ws.Range("N61").Value = -6676.8335
$ws.Range("H82").Value = 1283878.9
$ws.Range("I82").Value = 2349796.8
$ws.Range("J82").Value = 4777.4
$ws.Range("K82").Value = 2349796.8
$ws.Range("L82").Value = 4777.4
$ws.Range("M82").Value = -2349435.8
$ws.Range("N82").Value = -5499.4
$ws.Range("H85").Value = 1283878.9
$ws.Range("I85").Value = 2349796.8
$ws.Range("J85").Value = 4777.4
$ws.Range("K85").Value = 2349796.8
$ws.Range("L85").Value = 4777.4
$ws.Range("M85").Value = -2348548.8
$ws.Range("N85").Value = -7273.4
$ws.Range("H100").Value = 3270.818
$ws.Range("I100").Value = 3114.25
$ws.Range("J100").Value = 3360.2856
$ws.Range("K100").Value = 3114.25
$ws.Range("L100").Value = 3360.2856
$ws.Range("M100").Value = -2573.25
$ws.Range("N100").Value = -4442.2856
$ws.Range("H113").Value = 2634204
$ws.Range("I113").Value = 3847095.2
$ws.Range("J113").Value = 6272.8335
$ws.Range("K113").Value = 3847095.2
$ws.Range("L113").Value = 6272.8335
$ws.Range("M113").Value = -3844925.2
$ws.Range("N113").Value = -10612.8335
$ws.Range("H136").Value = 7600.6284
$ws.Range("I136").Value = 1710.1177
$ws.Range("J136").Value = 13163.889
$ws.Range("K136").Value = 5130.3531
$ws.Range("L136").Value = 39491.667
$ws.Range("M136").Value = -2580.3531
$ws.Range("N136").Value = -44591.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 4000
$ws.Range("I32").Value = 4000
$ws.Range("K32").Value = 4000
$ws.Range("M32").Value = -3683
$ws.Range("H81").Value = 8337314
$ws.Range("I81").Value = 1381.7142
$ws.Range("K81").Value = 2763.4284
$ws.Range("M81").Value = -1702.4284
$ws.Range("H82").Value = 23625
$ws.Range("I82").Value = 4250
$ws.Range("K82").Value = 4250
$ws.Range("M82").Value = -3867
$ws.Range("H84").Value = 8337314
$ws.Range("I84").Value = 1381.7142
$ws.Range("K84").Value = 13817.142
$ws.Range("M84").Value = -8513.142
$ws.Range("H85").Value = 23625
$ws.Range("I85").Value = 4250
$ws.Range("K85").Value = 4250
$ws.Range("M85").Value = -2924
$ws.Range("H107").Value = 1111.6364
$ws.Range("I107").Value = 1119.5714
$ws.Range("J107").Value = 1097.75
$ws.Range("K107").Value = 3358.7142
$ws.Range("L107").Value = 3293.25
$ws.Range("M107").Value = -1438.7142
$ws.Range("N107").Value = -7133.25
$ws.Range("H132").Value = 13893437
$ws.Range("I132").Value = 18523354
$ws.Range("K132").Value = 55570062
$ws.Range("M132").Value = -55567532
$ws.Range("H136").Value = 32295120
$ws.Range("I136").Value = 71429480
$ws.Range("J136").Value = 66823.82000000001
$ws.Range("K136").Value = 214288440
$ws.Range("L136").Value = 200471.46
$ws.Range("M136").Value = -214285890
$ws.Range("N136").Value = -205571.46
